$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.832.62'
$ws.Range('E2').Value = '  -0.04%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.360.51'
$ws.Range('E3').Value = '  +1.88%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.13'
$ws.Range('E5').Value = '  -0.45%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.10'
$ws.Range('E6').Value = '  -0.71%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.503'
$ws.Range('E7').Value = '  -0.77%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.484'
$ws.Range('E9').Value = '  -2.10%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.76'
$ws.Range('E10').Value = '  -1.87%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0784'
$ws.Range('E11').Value = '  -0.02%  '

$ws.Range('E12').Value = '  +2.44%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.18'
$ws.Range('E13').Value = '  -4.45%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.734.02'
$ws.Range('E14').Value = '  +2.17%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.70'
$ws.Range('E15').Value = '  -0.63%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.356.55'
$ws.Range('E16').Value = '  +1.70%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.794'
$ws.Range('E17').Value = '  +0.74%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.800.83'
$ws.Range('E18').Value = '  +0.05%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.96'
$ws.Range('E19').Value = '  -1.91%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.26'
$ws.Range('E20').Value = '  +1.58%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0883'
$ws.Range('E21').Value = '  -1.07%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.80'
$ws.Range('E22').Value = '  -0.02%  '

$ws.Range('E23').Value = '  -0.25%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.20'
$ws.Range('E24').Value = '  -2.20%  '

$ws.Range('E25').Value = '  -0.01%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.42'
$ws.Range('E26').Value = '  -0.32%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.58'
$ws.Range('E27').Value = '  +0.82%  '

$ws.Range('E28').Value = '  +0.28%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.20'
$ws.Range('E29').Value = '  +0.74%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.29'
$ws.Range('E30').Value = '  -2.95%  '

$ws.Range('E31').Value = '  +0.06%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.02'
$ws.Range('E32').Value = '  +0.25%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0724'
$ws.Range('E33').Value = '  +3.71%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.26'
$ws.Range('E34').Value = '  -3.61%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.85'
$ws.Range('E35').Value = '  +3.67%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.104'
$ws.Range('E36').Value = '  +4.05%  '

$ws.Range('E37').Value = '  -2.79%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.31'
$ws.Range('E38').Value = '  -1.18%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.77'
$ws.Range('E39').Value = '  +1.32%  '

$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '120.48'
$ws.Range('E40').Value = '  -27.63%  '

$ws.Range('E41').Value = '  -1.00%  '

$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.59'
$ws.Range('E42').Value = '  +3.92%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.930.27'
$ws.Range('E43').Value = '  +0.06%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0278'
$ws.Range('E44').Value = '  -0.20%  '

$ws.Range('E45').Value = '  +1.97%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.71'
$ws.Range('E46').Value = '  -1.83%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.14'
$ws.Range('E47').Value = '  -9.83%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.601.24'
$ws.Range('E48').Value = '  +2.24%  '

$ws.Range('E49').Value = '  +1.65%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '71.99'
$ws.Range('E50').Value = '  -0.37%  '

$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.14'
$ws.Range('E51').Value = '  +0.72%  '
